{"js": "// Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n// (percentages, dollar amounts, large numbers) inside specific resume bullet\n// paragraphs, matching the target OOXML diff.\n//\n// Strategy: for each target paragraph (identified by its exact, unique,\n// original plain text), scope a `search()` call to that paragraph only and\n// then bold + color each numeric metric substring found. Word's JS API\n// automatically splits the underlying <w:r> runs so that only the matched\n// text gets the new run properties, leaving the surrounding text in\n// separate (unformatted) runs - exactly the structure shown in the diff.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// List of [exact original paragraph text, [metric substrings to highlight, in order]]\nconst targets = [\n  [\n    \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    [\"23%\", \"64%\"],\n  ],\n  [\n    \"\\u2022 Utilized advanced sampling methods to decrease survey margin of error from \\u00B14.2% to \\u00B12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\",\n    [\"\\u00B14.2%\", \"\\u00B12.1%\", \"71%\", \"87%\"],\n  ],\n  [\n    \"\\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    [\"73.5%\", \"$4.7M\"],\n  ],\n  [\n    \"\\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n    [\"$2\"],\n  ],\n  [\n    \"\\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%\",\n    [\"57%\"],\n  ],\n  [\n    \"\\u2022 178% accuracy improvement in racial classification algorithms\",\n    [\"178%\"],\n  ],\n  [\n    \"\\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n    [\"73.5%\"],\n  ],\n  [\n    \"\\u2022 $4.7M savings enabled nonprofit access\",\n    [\"$4.7M\"],\n  ],\n  [\n    \"\\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\",\n    [\"12,847\"],\n  ],\n  [\n    \"\\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \\u00B14.2% to \\u00B12.1%\",\n    [\"\\u00B14.2%\", \"\\u00B12.1%\"],\n  ],\n  [\n    \"\\u2022 Increased voter turnout prediction accuracy from 71% to 87%\",\n    [\"71%\", \"87%\"],\n  ],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Build a lookup from exact paragraph text -> paragraph object.\nconst byText = new Map();\nfor (const p of paragraphs.items) {\n  byText.set(p.text, p);\n}\n\nfor (const [fullText, metrics] of targets) {\n  const paragraph = byText.get(fullText);\n  if (!paragraph) {\n    throw new Error(\"Could not locate target paragraph: \" + fullText);\n  }\n\n  // Find every metric substring, scoped to this paragraph only.\n  const ranges = [];\n  for (const metric of metrics) {\n    const found = paragraph.search(metric, { matchCase: true });\n    found.load(\"text\");\n    ranges.push(found);\n  }\n  await context.sync();\n\n  for (const found of ranges) {\n    if (found.items.length === 0) {\n      throw new Error(\"Metric not found in paragraph: \" + fullText);\n    }\n    const range = found.items[0];\n    range.font.bold = true;\n    range.font.color = HIGHLIGHT_COLOR;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n# (percentages, dollar amounts, large numbers) inside specific resume bullet\n# paragraphs, matching the target OOXML diff.\n#\n# Strategy: for each target paragraph (identified by its exact, unique,\n# original plain text, including the trailing paragraph-mark CR), use\n# Find.Execute scoped to that paragraph's Range to locate each numeric\n# metric substring in turn and apply Bold + the hybrid color (#2C3E50) to\n# just that substring. Word automatically splits the surrounding runs so\n# only the matched text receives the new run formatting - the same\n# structure produced by the diff.\n\n$d = $word.ActiveDocument\n\n$bullet = [char]0x2022\n$pm = [char]0x00B1\n\n# Hybrid highlight color 2C3E50 (RGB) expressed as a Word BGR long value.\n$highlightColor = 5258796\n\n# List of @(exact original paragraph text (no trailing CR), @(metric substrings, in order))\n$targets = @(\n    @(\n        \"$bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n        @(\"23%\", \"64%\")\n    ),\n    @(\n        \"$bullet Utilized advanced sampling methods to decrease survey margin of error from ${pm}4.2% to ${pm}2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\",\n        @(\"${pm}4.2%\", \"${pm}2.1%\", \"71%\", \"87%\")\n    ),\n    @(\n        \"$bullet Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\",\n        @(\"73.5%\", \"`$4.7M\")\n    ),\n    @(\n        \"$bullet Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\",\n        @(\"`$2\")\n    ),\n    @(\n        \"$bullet Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%\",\n        @(\"57%\")\n    ),\n    @(\n        \"$bullet 178% accuracy improvement in racial classification algorithms\",\n        @(\"178%\")\n    ),\n    @(\n        \"$bullet Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n        @(\"73.5%\")\n    ),\n    @(\n        \"$bullet `$4.7M savings enabled nonprofit access\",\n        @(\"`$4.7M\")\n    ),\n    @(\n        \"$bullet Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\",\n        @(\"12,847\")\n    ),\n    @(\n        \"$bullet Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ${pm}4.2% to ${pm}2.1%\",\n        @(\"${pm}4.2%\", \"${pm}2.1%\")\n    ),\n    @(\n        \"$bullet Increased voter turnout prediction accuracy from 71% to 87%\",\n        @(\"71%\", \"87%\")\n    )\n)\n\nforeach ($target in $targets) {\n    $fullText = $target[0]\n    $metrics = $target[1]\n\n    $paragraph = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -eq ($fullText + \"`r\")) {\n            $paragraph = $p\n            break\n        }\n    }\n    if ($paragraph -eq $null) {\n        throw \"Could not locate target paragraph: $fullText\"\n    }\n\n    # Use a range confined to the paragraph (excluding the end-of-paragraph\n    # mark) so Find cannot wander into neighboring paragraphs.\n    $searchRange = $paragraph.Range\n    $moveResult = $searchRange.MoveEnd(1, -1)\n\n    foreach ($metric in $metrics) {\n        $find = $searchRange.Find\n        $find.ClearFormatting()\n        $find.Text = $metric\n        $find.MatchCase = $true\n        $find.Forward = $true\n        $find.Wrap = 0\n        $found = $find.Execute()\n        if (-not $found) {\n            throw \"Metric '$metric' not found in paragraph: $fullText\"\n        }\n        $searchRange.Font.Bold = 1\n        $searchRange.Font.Color = $highlightColor\n    }\n}\n"}
